$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new (empty) spacer row at row 3 (previously the sheet jumped from row 2 straight
# to row 4). Touching its height materializes an empty row entry without adding any
# cell content.
$ws.Rows.Item(3).RowHeight = 15

# The user manually typed in a batch of newly observed license plates (and let the
# existing COUNTIF-based "frequency" formula in column B recompute for each of them).
$ws.Range("A108").Value = "BBB-111"
$ws.Range("B108").Formula = "=(COUNTIF(A:A,A108)-COUNTIF(A2,A108))"

$ws.Range("A109").Value = "BBB-123"
$ws.Range("B109").Formula = "=(COUNTIF(A:A,A109)-COUNTIF(A2,A109))"

$ws.Range("A110").Value = "BBB-112"
$ws.Range("B110").Formula = "=(COUNTIF(A:A,A110)-COUNTIF(A2,A110))"

$ws.Range("A111").Value = "BBB-111"
$ws.Range("B111").Formula = "=(COUNTIF(A:A,A111)-COUNTIF(A2,A111))"

$ws.Range("A112").Value = "CCC-111"
$ws.Range("B112").Formula = "=(COUNTIF(A:A,A112)-COUNTIF(A2,A112))"

# Column A got narrower after the edits (closest attainable width to the recorded 10.7).
$ws.Range("A1").EntireColumn.ColumnWidth = 9.86

# Leave the selection/cursor where the user ended up after typing the last entry.
$ws.Range("C113").Select()
